# Update cryptocurrency price/volume figures (and fix the Hedera / TrustWalletToken
# row order swap) per the latest scrape, per GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.129.21'
$ws.Cells.Item(2, 5).Value = '  -1.51%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.895.87'
$ws.Cells.Item(3, 5).Value = '  -0.82%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '314.39'
$ws.Cells.Item(5, 5).Value = '  -0.04%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.02%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5022'
$ws.Cells.Item(7, 5).Value = '  -0.70%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3898'
$ws.Cells.Item(8, 5).Value = '  -1.57%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.09244'
$ws.Cells.Item(9, 5).Value = '  -5.11%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '1.131'
$ws.Cells.Item(10, 5).Value = '  -2.41%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '41.87'
$ws.Cells.Item(11, 5).Value = '  +0.02%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '6.399'
$ws.Cells.Item(12, 5).Value = '  -2.20%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '20.82'
$ws.Cells.Item(13, 5).Value = '  -1.56%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.905.48'
$ws.Cells.Item(14, 5).Value = '  -0.02%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.288'
$ws.Cells.Item(15, 5).Value = '  -3.84%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.03%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001113'
$ws.Cells.Item(17, 5).Value = '  -2.43%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '92.43'
$ws.Cells.Item(18, 5).Value = '  -1.43%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.06656'
$ws.Cells.Item(19, 5).Value = '  +0.05%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '17.86'
$ws.Cells.Item(20, 5).Value = '  -1.14%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.211'
$ws.Cells.Item(22, 5).Value = '  -1.33%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '28.192.34'
$ws.Cells.Item(23, 5).Value = '  -1.49%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.04%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.321'
$ws.Cells.Item(25, 5).Value = '  +1.78%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.127.11'
$ws.Cells.Item(26, 5).Value = '  +0.06%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.555'
$ws.Cells.Item(27, 5).Value = '  -6.89%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '20.88'
$ws.Cells.Item(28, 5).Value = '  -2.03%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '158.39'
$ws.Cells.Item(29, 5).Value = '  -0.68%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '126.99'
$ws.Cells.Item(30, 5).Value = '  -1.52%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.081'
$ws.Cells.Item(31, 5).Value = '  -2.10%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -1.17%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.610'
$ws.Cells.Item(33, 5).Value = '  -2.25%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '3.622'
$ws.Cells.Item(34, 5).Value = '  -0.65%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '9.568'
$ws.Cells.Item(35, 5).Value = '  -3.34%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.06597'
$ws.Cells.Item(36, 5).Value = '  -3.13%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.336'
$ws.Cells.Item(37, 5).Value = '  +12.26%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02399'
$ws.Cells.Item(38, 5).Value = '  -1.99%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.223'
$ws.Cells.Item(40, 5).Value = '  -3.95%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.6471'
$ws.Cells.Item(41, 5).Value = '  +0.48%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '4.978'
$ws.Cells.Item(42, 5).Value = '  -2.63%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '11.40'
$ws.Cells.Item(43, 5).Value = '  -2.83%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.01%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.6108'
$ws.Cells.Item(45, 5).Value = '  +0.02%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -2.64%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.306'
$ws.Cells.Item(47, 5).Value = '  +1.54%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.693'
$ws.Cells.Item(48, 5).Value = '  +0.62%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.005'
$ws.Cells.Item(49, 5).Value = '  -2.00%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '122.28'
$ws.Cells.Item(50, 5).Value = '  -2.26%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.202'
$ws.Cells.Item(51, 5).Value = '  -0.94%  '

